$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 90914250
$ws.Range("I18").Value = 5605.5557
$ws.Range("K18").Value = 5605.5557
$ws.Range("M18").Value = -5321.5557

$ws.Range("H33").Value = 191.44444
$ws.Range("I33").Value = 191.44444
$ws.Range("K33").Value = 191.44444
$ws.Range("M33").Value = 37.55556000000001

$ws.Range("H113").Value = 5221.857
$ws.Range("I113").Value = 4767.6665
$ws.Range("K113").Value = 4767.6665
$ws.Range("M113").Value = -1513.6665

$ws.Range("H116").Value = 873617.6
$ws.Range("I116").Value = 997934.7
$ws.Range("K116").Value = 997934.7
$ws.Range("M116").Value = -994492.7

$ws.Range("H137").Value = 3807.3096
$ws.Range("I137").Value = 1683
$ws.Range("K137").Value = 5049
$ws.Range("M137").Value = -2499

$ws.Range("H141").Value = 4621.8887
$ws.Range("I141").Value = 4066.3333
$ws.Range("K141").Value = 12198.9999
$ws.Range("M141").Value = -7018.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 588.1786
$ws.Range("J2").Value = 650.2
$ws.Range("L2").Value = 650.2
$ws.Range("N2").Value = -876.2

$ws.Range("H45").Value = 46741.637
$ws.Range("I45").Value = 60646.625
$ws.Range("K45").Value = 60646.625
$ws.Range("M45").Value = -60269.625

$ws.Range("H61").Value = 4082.4583
$ws.Range("I61").Value = 2254.5715
$ws.Range("J61").Value = 6641.5
$ws.Range("K61").Value = 2254.5715
$ws.Range("L61").Value = 6641.5
$ws.Range("M61").Value = -2042.5715
$ws.Range("N61").Value = -7065.5

$ws.Range("H74").Value = 190712.27
$ws.Range("I74").Value = 371857.88
$ws.Range("J74").Value = 9566.666999999999
$ws.Range("K74").Value = 371857.88
$ws.Range("L74").Value = 9566.666999999999
$ws.Range("M74").Value = -370983.88
$ws.Range("N74").Value = -11314.667

$ws.Range("H77").Value = 190712.27
$ws.Range("I77").Value = 371857.88
$ws.Range("J77").Value = 9566.666999999999
$ws.Range("K77").Value = 1859289.4
$ws.Range("L77").Value = 47833.335
$ws.Range("M77").Value = -1854921.4
$ws.Range("N77").Value = -56569.335

$ws.Range("H102").Value = 4950.727
$ws.Range("I102").Value = 4800.8
$ws.Range("J102").Value = 6450
$ws.Range("K102").Value = 4800.8
$ws.Range("L102").Value = 6450
$ws.Range("M102").Value = -3178.8
$ws.Range("N102").Value = -9694

$ws.Range("H116").Value = 588.1786
$ws.Range("J116").Value = 650.2
$ws.Range("L116").Value = 650.2
$ws.Range("N116").Value = -5238.2

$ws.Range("H118").Value = 35000
$ws.Range("I118").Value = 35000
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 35000
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("M118").Value = -33343

$ws.Range("H132").Value = 2220.342
$ws.Range("I132").Value = 1645.7693
$ws.Range("J132").Value = 3465.25
$ws.Range("K132").Value = 4937.3079
$ws.Range("L132").Value = 10395.75
$ws.Range("M132").Value = -2407.3079
$ws.Range("N132").Value = -15455.75

$ws.Range("H136").Value = 4082.4583
$ws.Range("I136").Value = 2254.5715
$ws.Range("J136").Value = 6641.5
$ws.Range("K136").Value = 6763.7145
$ws.Range("L136").Value = 19924.5
$ws.Range("M136").Value = -4213.7145
$ws.Range("N136").Value = -25024.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 588.1786
$ws.Range("J3").Value = 650.2
$ws.Range("L3").Value = 650.2
$ws.Range("N3").Value = -878.2

$ws.Range("H86").Value = 3694.75
$ws.Range("I86").Value = 3529.6875
$ws.Range("J86").Value = 4024.875
$ws.Range("K86").Value = 3529.6875
$ws.Range("L86").Value = 4024.875
$ws.Range("M86").Value = -2406.6875
$ws.Range("N86").Value = -6270.875

$ws.Range("H89").Value = 3694.75
$ws.Range("I89").Value = 3529.6875
$ws.Range("J89").Value = 4024.875
$ws.Range("K89").Value = 17648.4375
$ws.Range("L89").Value = 20124.375
$ws.Range("M89").Value = -12032.4375
$ws.Range("N89").Value = -31356.375

$ws.Range("H94").Value = 38461990
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H134").Value = 2851.1904
$ws.Range("I134").Value = 2073.8462
$ws.Range("K134").Value = 6221.5386
$ws.Range("M134").Value = -3686.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3925.35
$ws.Range("I58").Value = 3138.7778
$ws.Range("K58").Value = 3138.7778
$ws.Range("M58").Value = -2935.7778

$ws.Range("H123").Value = 50000
$ws.Range("I123").Value = 50000
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 50000
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("M123").Value = -45100

$ws.Range("H132").Value = 11630685
$ws.Range("I132").Value = 13515823
$ws.Range("J132").Value = 5666.5
$ws.Range("K132").Value = 40547469
$ws.Range("L132").Value = 16999.5
$ws.Range("M132").Value = -40544939
$ws.Range("N132").Value = -22059.5

$ws.Range("H134").Value = 3233.7222
$ws.Range("I134").Value = 2880.4666
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 8641.399800000001
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -6106.399800000001
$ws.Range("N134").Value = -20070

$ws.Range("H136").Value = 3925.35
$ws.Range("I136").Value = 3138.7778
$ws.Range("K136").Value = 9416.3334
$ws.Range("M136").Value = -6866.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1909.5
$ws.Range("I70").Value = 379.33334
$ws.Range("K70").Value = 1138.00002
$ws.Range("M70").Value = -823.0000199999999

$ws.Range("H73").Value = 1909.5
$ws.Range("I73").Value = 379.33334
$ws.Range("K73").Value = 1138.00002
$ws.Range("M73").Value = -46.00001999999995

$ws.Range("H120").Value = 8000
$ws.Range("I120").Value = 8000
$ws.Range("K120").Value = 24000
$ws.Range("M120").Value = -19162

$ws.Range("H131").Value = 7199.25
$ws.Range("I131").Value = 14751.625
$ws.Range("J131").Value = 2164.3333
$ws.Range("K131").Value = 44254.875
$ws.Range("L131").Value = 6492.999899999999
$ws.Range("M131").Value = -39214.875
$ws.Range("N131").Value = -16572.9999

$ws.Range("H137").Value = 3439.3333
$ws.Range("I137").Value = 4072.6365
$ws.Range("J137").Value = 2742.7
$ws.Range("K137").Value = 12217.9095
$ws.Range("L137").Value = 8228.099999999999
$ws.Range("M137").Value = -7117.9095
$ws.Range("N137").Value = -18428.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1239
$ws.Range("I107").Value = 963.3333
$ws.Range("J107").Value = 1514.6666
$ws.Range("K107").Value = 963.3333
$ws.Range("L107").Value = 1514.6666
$ws.Range("M107").Value = 956.6667
$ws.Range("N107").Value = -5354.6666

$ws.Range("H116").Value = 44998
$ws.Range("J116").Value = 44998
$ws.Range("L116").Value = 44998
$ws.Range("N116").Value = -54176

$ws.Range("H122").Value = 2407961.8
$ws.Range("I122").Value = 3848459.5
$ws.Range("K122").Value = 11545378.5
$ws.Range("M122").Value = -11542928.5

$ws.Range("H132").Value = 3352.3076
$ws.Range("I132").Value = 3234.5454
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 9703.636200000001
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -7173.636200000001
$ws.Range("N132").Value = -17060

$ws.Range("H141").Value = 42784.715
$ws.Range("J141").Value = 42784.715
$ws.Range("L141").Value = 42784.715
$ws.Range("N141").Value = -53144.715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1118.5676
$ws.Range("I61").Value = 948.3226
$ws.Range("J61").Value = 1998.1666
$ws.Range("K61").Value = 948.3226
$ws.Range("L61").Value = 1998.1666
$ws.Range("M61").Value = -746.3226
$ws.Range("N61").Value = -2402.1666

$ws.Range("H68").Value = 7667.6665
$ws.Range("I68").Value = 7333.6665
$ws.Range("J68").Value = 8001.6665
$ws.Range("K68").Value = 7333.6665
$ws.Range("L68").Value = 8001.6665
$ws.Range("M68").Value = -6584.6665
$ws.Range("N68").Value = -9499.666499999999

$ws.Range("H71").Value = 7667.6665
$ws.Range("I71").Value = 7333.6665
$ws.Range("J71").Value = 8001.6665
$ws.Range("K71").Value = 36668.3325
$ws.Range("L71").Value = 40008.3325
$ws.Range("M71").Value = -32924.3325
$ws.Range("N71").Value = -47496.3325

$ws.Range("H113").Value = 1118.5676
$ws.Range("I113").Value = 948.3226
$ws.Range("J113").Value = 1998.1666
$ws.Range("K113").Value = 948.3226
$ws.Range("L113").Value = 1998.1666
$ws.Range("M113").Value = 1221.6774
$ws.Range("N113").Value = -6338.1666

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H132").Value = 2899.0247
$ws.Range("I132").Value = 2133.3242
$ws.Range("J132").Value = 10993.571
$ws.Range("K132").Value = 6399.9726
$ws.Range("L132").Value = 32980.713
$ws.Range("M132").Value = -3869.9726
$ws.Range("N132").Value = -38040.713

$ws.Range("H133").Value = 102319.25
$ws.Range("J133").Value = 102319.25
$ws.Range("L133").Value = 102319.25
$ws.Range("N133").Value = -107379.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6946900
$ws.Range("I132").Value = 8549640
$ws.Range("J132").Value = 1692.3334
$ws.Range("K132").Value = 25648920
$ws.Range("L132").Value = 5077.0002
$ws.Range("M132").Value = -25646390
$ws.Range("N132").Value = -10137.0002

Write-Output "Edit complete"